$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.744.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.81%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.006.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.94%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.29%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.998.50"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.85%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.496"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.10%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.66%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.65%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.09%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000223"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.42%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.484.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.99%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.747.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.71%  "

# Row 17
$ws.Range("E17").Value = "  -2.19%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.998.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.10%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.89%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "471.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.24%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.66%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.674"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.45%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.66%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.39%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.09%  "

# Row 26
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.65%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.30%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.27%  "

# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.69%  "

# Row 32
$ws.Range("B32").Value = "Mantle"
$ws.Range("C32").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.89%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.32%  "

# Row 34
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "55.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.05%  "

# Row 35
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.21%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "462.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.58%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.188.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.85%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0794"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.33%  "

# Row 40
$ws.Range("E40").Value = "  +1.14%  "

# Row 41
$ws.Range("E41").Value = "  -3.83%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.49%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.42%  "

# Row 44
$ws.Range("E44").Value = "  +0.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.29%  "

# Row 46
$ws.Range("E46").Value = "  -4.79%  "

# Row 47
$ws.Range("E47").Value = "  -3.20%  "

# Row 48
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.10%  "

# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.108"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.39%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0494"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.98%  "

# Row 51
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "
